# Fix for "out of bounds" bug in the branch-command table (sheet1, rows 35-39).
#
# The table previously mapped:
#   row35: bge  -> "if (t2 > t3) jump to label"      (wrong: that's bgt's description)
#   row36: blt  -> "if (t2 < t3) jump to label"       (wrong: that's blt but command was really for bge's slot)
#   row37: bltz -> "if (t2 < $zero) jump to label"
#   row38: (blank spacer row)
#   row39: label -> "some_label:"
#
# After the fix a missing "bgt" / ">=" comparison command is inserted, all the
# >,>=,<,<z rows are relabelled correctly, and the blank spacer + label rows
# shift down by one:
#   row35: bgt  -> "if (t2 > t3) jump to label"
#   row36: bge  -> "if (t2 >= t3) jump to label"      (new)
#   row37: blt  -> "if (t2 < t3) jump to label"
#   row38: bltz -> "if (t2 < $zero) jump to label"    (new row, was old row37)
#   row39: (blank spacer row, shifted down)
#   row40: label -> "some_label:" (shifted down)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Move-Row {
    param([int]$src, [int]$dst)
    # Copy the whole formatting of the source row onto the destination row first
    # (re-using existing style indices instead of letting a row/range Insert()
    # mint brand-new - and unused - style entries).
    $ws.Range("A$src`:F$src").Copy()
    $ws.Range("A$dst`:F$dst").PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = 0
    # Then copy the values over.
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($dst, $c).Value2 = $ws.Cells.Item($src, $c).Value2
    }
}

# Shift the blank spacer row (38) and the label row (39) down by one row to make
# room for the newly-inserted "bltz" row. Walk bottom-up so sources aren't
# clobbered before they are read.
Move-Row 39 40
Move-Row 38 39
Move-Row 37 38

# Row 35: commandName bge -> bgt ; description is unchanged ("if (t2 > t3) ...")
$ws.Range("A35").Value = "bgt"

# Row 36: commandName blt -> bge ; description becomes the new ">=" text
$ws.Range("A36").Value = "bge"
$ws.Range("B36").Value = "if (t2 >= t3) jump to label"

# Row 37: commandName bltz -> blt ; description becomes the existing "<" text
$ws.Range("A37").Value = "blt"
$ws.Range("B37").Value = "if (t2 < t3) jump to label"
# Column E (t3 placeholder) now participates, so it needs the same bold
# "60% - Accent5" style already used by D37/F37, reused via a formats-only paste
# (avoids minting a fresh, unused style entry).
$ws.Range("F37").Copy()
$ws.Range("E37").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("E37").Value = "t3"

# Update the view: scroll position and selection moved as part of the edit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C35").Select()
$excel.CutCopyMode = 0
